$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-17 16:45:08"

$wsZhCn.Range("H3").Value = "2016-08-17 16:44:57"
$wsZhCn.Range("K3").Value = "2016-08-17 16:45:30"

$wsDeDe.Range("H3").Value = "2016-08-17 16:45:08"
$wsDeDe.Range("K3").Value = "2016-08-17 16:45:39"
